$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 36:37 (pushes the old row 36 "ZZ/Page Clear" down to row 38)
$ws.Rows("36:37").Insert()
$ws.Rows("36:37").RowHeight = 14.25

# Copy the formatting pattern used by the other "categoria" data rows (e.g. row 30:
# column A bold/left category style, columns B:G plain text style) onto the two new rows
$ws.Range("A30:G30").Copy()
$ws.Range("A36:G37").PasteSpecial(-4122)

# Row 36: Rinascente Roma
$ws.Cells.Item(36, 1).Value = "Shopping Comercial"
$ws.Cells.Item(36, 2).Value = "Rinascente Roma"
$ws.Cells.Item(36, 3).Value = 41.904534018097301
$ws.Cells.Item(36, 4).Value = 12.4846873410065
$ws.Cells.Item(36, 5).Value = "bolsas.png"
$ws.Cells.Item(36, 6).Value = "Via del Tritone, 61, 00187 Roma RM, Itália"
$ws.Cells.Item(36, 7).Value = "Horário: 10:00–21:30"

# Row 37: Galleria Alberto Sordi
$ws.Cells.Item(37, 1).Value = "Shopping Comercial"
$ws.Cells.Item(37, 2).Value = "Galleria Alberto Sordi"
$ws.Cells.Item(37, 3).Value = 41.903640774176502
$ws.Cells.Item(37, 4).Value = 12.4815967890303
$ws.Cells.Item(37, 5).Value = "bolsas.png"
$ws.Cells.Item(37, 6).Value = "P.za Colonna, 00187 Roma RM, Itália"
$ws.Cells.Item(37, 7).Value = "Horário: 09:00–20:00"

# Match the workbook's saved view/selection state
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G38").Select()
